$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8858482241630554
$ws.Range("B1").Value = 1.77877938747406
$ws.Range("D1").Value = 1.889378309249878
$ws.Range("E1").Value = 1.118734359741211
